$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.123.90"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').Value = "'2.050.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.97%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'231.89"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').Value = "'0.619"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.52%  '
$ws.Range('D8').Value = "'56.92"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.20%  '
$ws.Range('D9').Value = "'0.382"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.65%  '
$ws.Range('D10').Value = "'57.57"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.63%  '
$ws.Range('D11').Value = "'0.0756"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = "'2.354.83"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = "'14.41"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.42%  '
$ws.Range('D15').Value = "'20.75"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.53%  '
$ws.Range('D16').Value = "'0.774"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').Value = "'5.12"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').Value = "'2.049.79"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = "'37.095.84"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').Value = "'6.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.28%  '
$ws.Range('D21').Value = "'69.21"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.91%  '
$ws.Range('E22').Value = '  +1.54%  '
$ws.Range('D23').Value = "'224.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').Value = "'165.55"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('D28').Value = "'8.75"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.00%  '
$ws.Range('E29').Value = '  +6.59%  '
$ws.Range('D30').Value = "'18.96"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').Value = "'0.126"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').Value = "'4.44"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('D34').Value = "'0.0615"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('D35').Value = "'4.59"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.63%  '
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').Value = "'1.478.17"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = "'96.46"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.10%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').Value = "'4.40"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('D45').Value = "'1.17"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.63%  '
$ws.Range('D46').Value = "'0.0929"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('D47').Value = "'0.0210"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.52%  '
$ws.Range('D48').Value = "'1.02"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.95%  '
$ws.Range('D49').Value = "'7.16"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.68%  '
$ws.Range('D50').Value = "'15.05"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').Value = "'2.94"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.58%  '
